$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("correl-employment-GES")

# Update the correlation commentary cell to the refreshed analysis text
$ws.Range("C10").Value = "There is a very strong positive correlation between Overall Employment and Graduate Employment"

# Row 10 grew taller to fit the updated wrapped text
$ws.Rows.Item(10).RowHeight = 100.8

# Reflect the author's final selection on the sheet
[void]$ws.Range("O10").Select()

# Drop the stale/orphaned chart cache names left over from earlier chart edits
for ($i = 0; $i -le 6; $i++) {
    $name = "_xlchart.v1.$i"
    [void]$wb.Names.Item($name).Delete()
}
